$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): Right count 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): Right count 63 -> 105
$ws.Range("B12").Value = 105

# Update the Max (score fraction) text: "60/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
